$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status text for the second data row (Status columns) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
# 39.15 is the ColumnWidth input that rounds to a stored XML width of exactly 40
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15
$wsZhCn.Range("P3").Value = "Handback file name: yfqgjjyt.jbf is different with handoff file name: fe4a299e-294c-4690-8970-fbc2159bca0b.6a864893672f3ac796fda4ca581ab793aa0278ea.zh-cn."

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
$wsDeDe.Range("P3").Value = "Handback file name: yfqgjjyt.jbf is different with handoff file name: fe4a299e-294c-4690-8970-fbc2159bca0b.6a864893672f3ac796fda4ca581ab793aa0278ea.de-de."
